$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object "object[,]" 72,4
$data[0,0] = "Combined Ratio"
$data[0,1] = 77.09999999999999
$data[0,2] = "Intact"
$data[0,3] = "Q4 2024"
$data[1,0] = "Combined Ratio"
$data[1,1] = 75.8
$data[1,2] = "Intact"
$data[1,3] = "Q4 2023"
$data[2,0] = "Claims Ratio"
$data[2,1] = 43.3
$data[2,2] = "Intact"
$data[2,3] = "Q4 2023"
$data[3,0] = "Claims Ratio"
$data[3,1] = 42.6
$data[3,2] = "Intact"
$data[3,3] = "Q4 2024"
$data[4,0] = "Core Claim Ratio"
$data[4,1] = 41.6
$data[4,2] = "Intact"
$data[4,3] = "Q4 2024"
$data[5,0] = "Core Claim Ratio"
$data[5,1] = 44.4
$data[5,2] = "Intact"
$data[5,3] = "Q4 2023"
$data[6,0] = "CAT Loss Ratio"
$data[6,1] = 0.6
$data[6,2] = "Intact"
$data[6,3] = "Q4 2023"
$data[7,0] = "CAT Loss Ratio"
$data[7,1] = 3.9
$data[7,2] = "Intact"
$data[7,3] = "Q4 2024"
$data[8,0] = "Expense Ratio"
$data[8,1] = 34.5
$data[8,2] = "Intact"
$data[8,3] = "Q4 2024"
$data[9,0] = "Expense Ratio"
$data[9,1] = 32.5
$data[9,2] = "Intact"
$data[9,3] = "Q4 2023"
$data[10,0] = "PYD Ratio"
$data[10,1] = -1.7
$data[10,2] = "Intact"
$data[10,3] = "Q4 2023"
$data[11,0] = "PYD Ratio"
$data[11,1] = -2.9
$data[11,2] = "Intact"
$data[11,3] = "Q4 2024"
$data[12,0] = "Gross Written Premium"
$data[12,1] = 1030
$data[12,2] = "Intact"
$data[12,3] = "Q4 2024"
$data[13,0] = "Gross Written Premium"
$data[13,1] = 946
$data[13,2] = "Intact"
$data[13,3] = "Q4 2023"
$data[14,0] = "Underwriting Income"
$data[14,1] = 229
$data[14,2] = "Intact"
$data[14,3] = "Q4 2023"
$data[15,0] = "Underwriting Income"
$data[15,1] = 237
$data[15,2] = "Intact"
$data[15,3] = "Q4 2024"
$data[16,0] = "ROE"
$data[16,1] = 14.2
$data[16,2] = "Intact"
$data[16,3] = "Q4 2024"
$data[17,0] = "ROE"
$data[17,1] = 8.800000000000001
$data[17,2] = "Intact"
$data[17,3] = "Q4 2023"
$data[18,0] = "ROE"
$data[18,1] = 8.800000000000001
$data[18,2] = "Intact"
$data[18,3] = "2023"
$data[19,0] = "ROE"
$data[19,1] = 14.2
$data[19,2] = "Intact"
$data[19,3] = "2024"
$data[20,0] = "Combined Ratio"
$data[20,1] = 96.5
$data[20,2] = "Intact"
$data[20,3] = "2024"
$data[21,0] = "Combined Ratio"
$data[21,1] = 100.7
$data[21,2] = "Intact"
$data[21,3] = "2023"
$data[22,0] = "Claims Ratio"
$data[22,1] = 67
$data[22,2] = "Intact"
$data[22,3] = "2023"
$data[23,0] = "Claims Ratio"
$data[23,1] = 62
$data[23,2] = "Intact"
$data[23,3] = "2024"
$data[24,0] = "Core Claim Ratio"
$data[24,1] = 46.1
$data[24,2] = "Intact"
$data[24,3] = "2024"
$data[25,0] = "Core Claim Ratio"
$data[25,1] = 49
$data[25,2] = "Intact"
$data[25,3] = "2023"
$data[26,0] = "CAT Loss Ratio"
$data[26,1] = 18.3
$data[26,2] = "Intact"
$data[26,3] = "2023"
$data[27,0] = "CAT Loss Ratio"
$data[27,1] = 19.7
$data[27,2] = "Intact"
$data[27,3] = "2024"
$data[28,0] = "Expense Ratio"
$data[28,1] = 34.5
$data[28,2] = "Intact"
$data[28,3] = "2024"
$data[29,0] = "Expense Ratio"
$data[29,1] = 33.7
$data[29,2] = "Intact"
$data[29,3] = "2023"
$data[30,0] = "PYD Ratio"
$data[30,1] = -0.3
$data[30,2] = "Intact"
$data[30,3] = "2023"
$data[31,0] = "PYD Ratio"
$data[31,1] = -3.8
$data[31,2] = "Intact"
$data[31,3] = "2024"
$data[32,0] = "Gross Written Premium"
$data[32,1] = 4222
$data[32,2] = "Intact"
$data[32,3] = "2024"
$data[33,0] = "Gross Written Premium"
$data[33,1] = 3877
$data[33,2] = "Intact"
$data[33,3] = "2023"
$data[34,0] = "Underwriting Income"
$data[34,1] = -26
$data[34,2] = "Intact"
$data[34,3] = "2023"
$data[35,0] = "Underwriting Income"
$data[35,1] = 138
$data[35,2] = "Intact"
$data[35,3] = "2024"
$data[36,0] = "Combined Ratio"
$data[36,1] = 96.3
$data[36,2] = "Definity"
$data[36,3] = "2024"
$data[37,0] = "Combined Ratio"
$data[37,1] = 99.3
$data[37,2] = "Definity"
$data[37,3] = "2023"
$data[38,0] = "Claims Ratio"
$data[38,1] = 64.5
$data[38,2] = "Definity"
$data[38,3] = "2023"
$data[39,0] = "Claims Ratio"
$data[39,1] = 62.7
$data[39,2] = "Definity"
$data[39,3] = "2024"
$data[40,0] = "Core Claim Ratio"
$data[40,1] = 49.9
$data[40,2] = "Definity"
$data[40,3] = "2024"
$data[41,0] = "Core Claim Ratio"
$data[41,1] = 50
$data[41,2] = "Definity"
$data[41,3] = "2023"
$data[42,0] = "CAT Loss Ratio"
$data[42,1] = 15.6
$data[42,2] = "Definity"
$data[42,3] = "2023"
$data[43,0] = "CAT Loss Ratio"
$data[43,1] = 15.7
$data[43,2] = "Definity"
$data[43,3] = "2024"
$data[44,0] = "Expense Ratio"
$data[44,1] = 33.6
$data[44,2] = "Definity"
$data[44,3] = "2024"
$data[45,0] = "Expense Ratio"
$data[45,1] = 34.8
$data[45,2] = "Definity"
$data[45,3] = "2023"
$data[46,0] = "PYD Ratio"
$data[46,1] = -1.1
$data[46,2] = "Definity"
$data[46,3] = "2023"
$data[47,0] = "PYD Ratio"
$data[47,1] = -2.9
$data[47,2] = "Definity"
$data[47,3] = "2024"
$data[48,0] = "Gross Written Premium"
$data[48,1] = 1183.9
$data[48,2] = "Definity"
$data[48,3] = "2024"
$data[49,0] = "Gross Written Premium"
$data[49,1] = 1113.1
$data[49,2] = "Definity"
$data[49,3] = "2023"
$data[50,0] = "Underwriting Income"
$data[50,1] = 7.2
$data[50,2] = "Definity"
$data[50,3] = "2023"
$data[51,0] = "Underwriting Income"
$data[51,1] = 40.1
$data[51,2] = "Definity"
$data[51,3] = "2024"
$data[52,0] = "ROE"
$data[52,1] = 10.6
$data[52,2] = "Definity"
$data[52,3] = "2024"
$data[53,0] = "ROE"
$data[53,1] = 9.199999999999999
$data[53,2] = "Definity"
$data[53,3] = "2023"
$data[54,0] = "Combined Ratio"
$data[54,1] = 80.09999999999999
$data[54,2] = "Definity"
$data[54,3] = "Q4 2023"
$data[55,0] = "Combined Ratio"
$data[55,1] = 82.8
$data[55,2] = "Definity"
$data[55,3] = "Q4 2024"
$data[56,0] = "Claims Ratio"
$data[56,1] = 50.4
$data[56,2] = "Definity"
$data[56,3] = "Q4 2024"
$data[57,0] = "Claims Ratio"
$data[57,1] = 46.6
$data[57,2] = "Definity"
$data[57,3] = "Q4 2023"
$data[58,0] = "Core Claim Ratio"
$data[58,1] = 47
$data[58,2] = "Definity"
$data[58,3] = "Q4 2023"
$data[59,0] = "Core Claim Ratio"
$data[59,1] = 46.2
$data[59,2] = "Definity"
$data[59,3] = "Q4 2024"
$data[60,0] = "CAT Loss Ratio"
$data[60,1] = 7.2
$data[60,2] = "Definity"
$data[60,3] = "Q4 2024"
$data[61,0] = "CAT Loss Ratio"
$data[61,1] = 1.2
$data[61,2] = "Definity"
$data[61,3] = "Q4 2023"
$data[62,0] = "Expense Ratio"
$data[62,1] = 33.5
$data[62,2] = "Definity"
$data[62,3] = "Q4 2023"
$data[63,0] = "Expense Ratio"
$data[63,1] = 32.4
$data[63,2] = "Definity"
$data[63,3] = "Q4 2024"
$data[64,0] = "PYD Ratio"
$data[64,1] = -3
$data[64,2] = "Definity"
$data[64,3] = "Q4 2024"
$data[65,0] = "PYD Ratio"
$data[65,1] = -1.6
$data[65,2] = "Definity"
$data[65,3] = "Q4 2023"
$data[66,0] = "Gross Written Premium"
$data[66,1] = 278
$data[66,2] = "Definity"
$data[66,3] = "Q4 2023"
$data[67,0] = "Gross Written Premium"
$data[67,1] = 295.2
$data[67,2] = "Definity"
$data[67,3] = "Q4 2024"
$data[68,0] = "Underwriting Income"
$data[68,1] = 47.9
$data[68,2] = "Definity"
$data[68,3] = "Q4 2024"
$data[69,0] = "Underwriting Income"
$data[69,1] = 53.1
$data[69,2] = "Definity"
$data[69,3] = "Q4 2023"
$data[70,0] = "ROE"
$data[70,1] = 9.199999999999999
$data[70,2] = "Definity"
$data[70,3] = "Q4 2023"
$data[71,0] = "ROE"
$data[71,1] = 10.6
$data[71,2] = "Definity"
$data[71,3] = "Q4 2024"
$range = $ws.Range("A130:D201")
$range.Value = $data
